$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 71.41194166666666
$ws.Cells.Item(2, 8).Value = 214.235825
$ws.Cells.Item(2, 9).Value = 0.02299241149786563
$ws.Cells.Item(2, 10).Value = 0.02299241149786563
$ws.Cells.Item(2, 13).Value = 0.7592313333333333
$ws.Cells.Item(2, 14).Value = 2.277694
$ws.Cells.Item(2, 15).Value = 0.05311309006881704
$ws.Cells.Item(2, 16).Value = 0.05311309006881704
$ws.Cells.Item(2, 17).Value = 54.21818368750555
$ws.Cells.Item(2, 18).Value = 487.96365318755
$ws.Cells.Item(2, 19).Value = 0.001221198022785442
$ws.Cells.Item(2, 20).Value = 0.001221198022785442

# Row 3
$ws.Cells.Item(3, 7).Value = 71.41194166666666
$ws.Cells.Item(3, 8).Value = 214.235825
$ws.Cells.Item(3, 9).Value = 0.02299241149786563
$ws.Cells.Item(3, 10).Value = 0.02299241149786563
$ws.Cells.Item(3, 15).Value = 0.1721303313829956
$ws.Cells.Item(3, 16).Value = 0.1721303313829956
$ws.Cells.Item(3, 17).Value = 175.7117485166555
$ws.Cells.Item(3, 18).Value = 1581.4057366499
$ws.Cells.Item(3, 19).Value = 0.003957691410421809
$ws.Cells.Item(3, 20).Value = 0.00395769141042181

# Row 4
$ws.Cells.Item(4, 7).Value = 71.41194166666666
$ws.Cells.Item(4, 8).Value = 214.235825
$ws.Cells.Item(4, 9).Value = 0.02299241149786563
$ws.Cells.Item(4, 10).Value = 0.02299241149786563
$ws.Cells.Item(4, 13).Value = 1.018760333333333
$ws.Cells.Item(4, 14).Value = 3.056281
$ws.Cells.Item(4, 15).Value = 0.07126880433834142
$ws.Cells.Item(4, 16).Value = 0.07126880433834142
$ws.Cells.Item(4, 17).Value = 72.75165349631388
$ws.Cells.Item(4, 18).Value = 654.7648814668249
$ws.Cells.Item(4, 19).Value = 0.001638641676308017
$ws.Cells.Item(4, 20).Value = 0.001638641676308017

# Row 5
$ws.Cells.Item(5, 7).Value = 71.41194166666666
$ws.Cells.Item(5, 8).Value = 214.235825
$ws.Cells.Item(5, 9).Value = 0.02299241149786563
$ws.Cells.Item(5, 10).Value = 0.02299241149786563
$ws.Cells.Item(5, 13).Value = 9.893154666666666
$ws.Cells.Item(5, 14).Value = 29.679464
$ws.Cells.Item(5, 15).Value = 0.6920894749804902
$ws.Cells.Item(5, 16).Value = 0.6920894749804903
$ws.Cells.Item(5, 17).Value = 706.489383955311
$ws.Cells.Item(5, 18).Value = 6358.4044555978
$ws.Cells.Item(5, 19).Value = 0.01591280600209321
$ws.Cells.Item(5, 20).Value = 0.01591280600209322

# Row 6
$ws.Cells.Item(6, 7).Value = 71.41194166666666
$ws.Cells.Item(6, 8).Value = 214.235825
$ws.Cells.Item(6, 9).Value = 0.02299241149786563
$ws.Cells.Item(6, 10).Value = 0.02299241149786563
$ws.Cells.Item(6, 13).Value = 0.1629343333333333
$ws.Cells.Item(6, 14).Value = 0.488803
$ws.Cells.Item(6, 15).Value = 0.01139829922935564
$ws.Cells.Item(6, 16).Value = 0.01139829922935564
$ws.Cells.Item(6, 17).Value = 11.63545710749722
$ws.Cells.Item(6, 18).Value = 104.719113967475
$ws.Cells.Item(6, 19).Value = 0.0002620743862571497
$ws.Cells.Item(6, 20).Value = 0.0002620743862571497

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9549836193138445
$ws.Cells.Item(7, 10).Value = 0.9549836193138445
$ws.Cells.Item(7, 13).Value = 0.7592313333333333
$ws.Cells.Item(7, 14).Value = 2.277694
$ws.Cells.Item(7, 15).Value = 0.05311309006881704
$ws.Cells.Item(7, 16).Value = 0.05311309006881704
$ws.Cells.Item(7, 17).Value = 2251.937657575559
$ws.Cells.Item(7, 18).Value = 20267.43891818004
$ws.Cells.Item(7, 19).Value = 0.0507221309868611
$ws.Cells.Item(7, 20).Value = 0.0507221309868611

# Row 8
$ws.Cells.Item(8, 9).Value = 0.9549836193138445
$ws.Cells.Item(8, 10).Value = 0.9549836193138445
$ws.Cells.Item(8, 15).Value = 0.1721303313829956
$ws.Cells.Item(8, 16).Value = 0.1721303313829956
$ws.Cells.Item(8, 19).Value = 0.1643816468578246
$ws.Cells.Item(8, 20).Value = 0.1643816468578246

# Row 9
$ws.Cells.Item(9, 9).Value = 0.9549836193138445
$ws.Cells.Item(9, 10).Value = 0.9549836193138445
$ws.Cells.Item(9, 13).Value = 1.018760333333333
$ws.Cells.Item(9, 14).Value = 3.056281
$ws.Cells.Item(9, 15).Value = 0.07126880433834142
$ws.Cells.Item(9, 16).Value = 0.07126880433834142
$ws.Cells.Item(9, 17).Value = 3021.720334703735
$ws.Cells.Item(9, 18).Value = 27195.48301233361
$ws.Cells.Item(9, 19).Value = 0.0680605407111995
$ws.Cells.Item(9, 20).Value = 0.0680605407111995

# Row 10
$ws.Cells.Item(10, 9).Value = 0.9549836193138445
$ws.Cells.Item(10, 10).Value = 0.9549836193138445
$ws.Cells.Item(10, 13).Value = 9.893154666666666
$ws.Cells.Item(10, 14).Value = 29.679464
$ws.Cells.Item(10, 15).Value = 0.6920894749804902
$ws.Cells.Item(10, 16).Value = 0.6920894749804903
$ws.Cells.Item(10, 17).Value = 29343.84629289893
$ws.Cells.Item(10, 18).Value = 264094.6166360904
$ws.Cells.Item(10, 19).Value = 0.660934111705887
$ws.Cells.Item(10, 20).Value = 0.6609341117058871

# Row 11
$ws.Cells.Item(11, 9).Value = 0.9549836193138445
$ws.Cells.Item(11, 10).Value = 0.9549836193138445
$ws.Cells.Item(11, 13).Value = 0.1629343333333333
$ws.Cells.Item(11, 14).Value = 0.488803
$ws.Cells.Item(11, 15).Value = 0.01139829922935564
$ws.Cells.Item(11, 16).Value = 0.01139829922935564
$ws.Cells.Item(11, 17).Value = 483.2755773321202
$ws.Cells.Item(11, 18).Value = 4349.480195989082
$ws.Cells.Item(11, 19).Value = 0.01088518905207226
$ws.Cells.Item(11, 20).Value = 0.01088518905207226

# Row 12
$ws.Cells.Item(12, 7).Value = 1.469787333333333
$ws.Cells.Item(12, 8).Value = 4.409362
$ws.Cells.Item(12, 9).Value = 0.0004732255473474234
$ws.Cells.Item(12, 10).Value = 0.0004732255473474233
$ws.Cells.Item(12, 13).Value = 0.7592313333333333
$ws.Cells.Item(12, 14).Value = 2.277694
$ws.Cells.Item(12, 15).Value = 0.05311309006881704
$ws.Cells.Item(12, 16).Value = 0.05311309006881704
$ws.Cells.Item(12, 17).Value = 1.115908596803111
$ws.Cells.Item(12, 18).Value = 10.043177371228
$ws.Cells.Item(12, 19).Value = 0.00002513447111912894
$ws.Cells.Item(12, 20).Value = 0.00002513447111912893

# Row 13
$ws.Cells.Item(13, 7).Value = 1.469787333333333
$ws.Cells.Item(13, 8).Value = 4.409362
$ws.Cells.Item(13, 9).Value = 0.0004732255473474234
$ws.Cells.Item(13, 10).Value = 0.0004732255473474233
$ws.Cells.Item(13, 15).Value = 0.1721303313829956
$ws.Cells.Item(13, 16).Value = 0.1721303313829956
$ws.Cells.Item(13, 17).Value = 3.616466605727111
$ws.Cells.Item(13, 18).Value = 32.548199451544
$ws.Cells.Item(13, 19).Value = 0.00008145647028381145
$ws.Cells.Item(13, 20).Value = 0.00008145647028381145

# Row 14
$ws.Cells.Item(14, 7).Value = 1.469787333333333
$ws.Cells.Item(14, 8).Value = 4.409362
$ws.Cells.Item(14, 9).Value = 0.0004732255473474234
$ws.Cells.Item(14, 10).Value = 0.0004732255473474233
$ws.Cells.Item(14, 13).Value = 1.018760333333333
$ws.Cells.Item(14, 14).Value = 3.056281
$ws.Cells.Item(14, 15).Value = 0.07126880433834142
$ws.Cells.Item(14, 16).Value = 0.07126880433834142
$ws.Cells.Item(14, 17).Value = 1.497361033635778
$ws.Cells.Item(14, 18).Value = 13.476249302722
$ws.Cells.Item(14, 19).Value = 0.00003372621894180804
$ws.Cells.Item(14, 20).Value = 0.00003372621894180803

# Row 15
$ws.Cells.Item(15, 7).Value = 1.469787333333333
$ws.Cells.Item(15, 8).Value = 4.409362
$ws.Cells.Item(15, 9).Value = 0.0004732255473474234
$ws.Cells.Item(15, 10).Value = 0.0004732255473474233
$ws.Cells.Item(15, 13).Value = 9.893154666666666
$ws.Cells.Item(15, 14).Value = 29.679464
$ws.Cells.Item(15, 15).Value = 0.6920894749804902
$ws.Cells.Item(15, 16).Value = 0.6920894749804903
$ws.Cells.Item(15, 17).Value = 14.54083341577422
$ws.Cells.Item(15, 18).Value = 130.867500741968
$ws.Cells.Item(15, 19).Value = 0.0003275144206110334
$ws.Cells.Item(15, 20).Value = 0.0003275144206110334

# Row 16
$ws.Cells.Item(16, 7).Value = 1.469787333333333
$ws.Cells.Item(16, 8).Value = 4.409362
$ws.Cells.Item(16, 9).Value = 0.0004732255473474234
$ws.Cells.Item(16, 10).Value = 0.0004732255473474233
$ws.Cells.Item(16, 13).Value = 0.1629343333333333
$ws.Cells.Item(16, 14).Value = 0.488803
$ws.Cells.Item(16, 15).Value = 0.01139829922935564
$ws.Cells.Item(16, 16).Value = 0.01139829922935564
$ws.Cells.Item(16, 17).Value = 0.2394788192984444
$ws.Cells.Item(16, 18).Value = 2.155309373686
$ws.Cells.Item(16, 19).Value = 0.000005393966391641539
$ws.Cells.Item(16, 20).Value = 0.000005393966391641539

# Row 17
$ws.Cells.Item(17, 7).Value = 65.51927933333333
$ws.Cells.Item(17, 8).Value = 196.557838
$ws.Cells.Item(17, 9).Value = 0.02109515854515373
$ws.Cells.Item(17, 10).Value = 0.02109515854515373
$ws.Cells.Item(17, 13).Value = 0.7592313333333333
$ws.Cells.Item(17, 14).Value = 2.277694
$ws.Cells.Item(17, 15).Value = 0.05311309006881704
$ws.Cells.Item(17, 16).Value = 0.05311309006881704
$ws.Cells.Item(17, 17).Value = 49.74428980728577
$ws.Cells.Item(17, 18).Value = 447.6986082655719
$ws.Cells.Item(17, 19).Value = 0.001120429055824726
$ws.Cells.Item(17, 20).Value = 0.001120429055824726

# Row 18
$ws.Cells.Item(18, 7).Value = 65.51927933333333
$ws.Cells.Item(18, 8).Value = 196.557838
$ws.Cells.Item(18, 9).Value = 0.02109515854515373
$ws.Cells.Item(18, 10).Value = 0.02109515854515373
$ws.Cells.Item(18, 15).Value = 0.1721303313829956
$ws.Cells.Item(18, 16).Value = 0.1721303313829956
$ws.Cells.Item(18, 17).Value = 161.2126328527618
$ws.Cells.Item(18, 18).Value = 1450.913695674856
$ws.Cells.Item(18, 19).Value = 0.003631116630954143
$ws.Cells.Item(18, 20).Value = 0.003631116630954143

# Row 19
$ws.Cells.Item(19, 7).Value = 65.51927933333333
$ws.Cells.Item(19, 8).Value = 196.557838
$ws.Cells.Item(19, 9).Value = 0.02109515854515373
$ws.Cells.Item(19, 10).Value = 0.02109515854515373
$ws.Cells.Item(19, 13).Value = 1.018760333333333
$ws.Cells.Item(19, 14).Value = 3.056281
$ws.Cells.Item(19, 15).Value = 0.07126880433834142
$ws.Cells.Item(19, 16).Value = 0.07126880433834142
$ws.Cells.Item(19, 17).Value = 66.74844285338644
$ws.Cells.Item(19, 18).Value = 600.7359856804779
$ws.Cells.Item(19, 19).Value = 0.001503426726840852
$ws.Cells.Item(19, 20).Value = 0.001503426726840852

# Row 20
$ws.Cells.Item(20, 7).Value = 65.51927933333333
$ws.Cells.Item(20, 8).Value = 196.557838
$ws.Cells.Item(20, 9).Value = 0.02109515854515373
$ws.Cells.Item(20, 10).Value = 0.02109515854515373
$ws.Cells.Item(20, 13).Value = 9.893154666666666
$ws.Cells.Item(20, 14).Value = 29.679464
$ws.Cells.Item(20, 15).Value = 0.6920894749804902
$ws.Cells.Item(20, 16).Value = 0.6920894749804903
$ws.Cells.Item(20, 17).Value = 648.1923640932034
$ws.Cells.Item(20, 18).Value = 5833.731276838831
$ws.Cells.Item(20, 19).Value = 0.01459973720214565
$ws.Cells.Item(20, 20).Value = 0.01459973720214565

# Row 21
$ws.Cells.Item(21, 7).Value = 65.51927933333333
$ws.Cells.Item(21, 8).Value = 196.557838
$ws.Cells.Item(21, 9).Value = 0.02109515854515373
$ws.Cells.Item(21, 10).Value = 0.02109515854515373
$ws.Cells.Item(21, 13).Value = 0.1629343333333333
$ws.Cells.Item(21, 14).Value = 0.488803
$ws.Cells.Item(21, 15).Value = 0.01139829922935564
$ws.Cells.Item(21, 16).Value = 0.01139829922935564
$ws.Cells.Item(21, 17).Value = 10.67534009865711
$ws.Cells.Item(21, 18).Value = 96.07806088791399
$ws.Cells.Item(21, 19).Value = 0.0002404489293883609
$ws.Cells.Item(21, 20).Value = 0.0002404489293883609

# Row 22
$ws.Cells.Item(22, 7).Value = 1.414998
$ws.Cells.Item(22, 8).Value = 4.244994
$ws.Cells.Item(22, 9).Value = 0.0004555850957885808
$ws.Cells.Item(22, 10).Value = 0.0004555850957885808
$ws.Cells.Item(22, 13).Value = 0.7592313333333333
$ws.Cells.Item(22, 14).Value = 2.277694
$ws.Cells.Item(22, 15).Value = 0.05311309006881704
$ws.Cells.Item(22, 16).Value = 0.05311309006881704
$ws.Cells.Item(22, 17).Value = 1.074310818204
$ws.Cells.Item(22, 18).Value = 9.668797363835999
$ws.Cells.Item(22, 19).Value = 0.00002419753222662953
$ws.Cells.Item(22, 20).Value = 0.00002419753222662953

# Row 23
$ws.Cells.Item(23, 7).Value = 1.414998
$ws.Cells.Item(23, 8).Value = 4.244994
$ws.Cells.Item(23, 9).Value = 0.0004555850957885808
$ws.Cells.Item(23, 10).Value = 0.0004555850957885808
$ws.Cells.Item(23, 15).Value = 0.1721303313829956
$ws.Cells.Item(23, 16).Value = 0.1721303313829956
$ws.Cells.Item(23, 17).Value = 3.481655405592
$ws.Cells.Item(23, 18).Value = 31.334898650328
$ws.Cells.Item(23, 19).Value = 0.0000784200135112422
$ws.Cells.Item(23, 20).Value = 0.0000784200135112422

# Row 24
$ws.Cells.Item(24, 7).Value = 1.414998
$ws.Cells.Item(24, 8).Value = 4.244994
$ws.Cells.Item(24, 9).Value = 0.0004555850957885808
$ws.Cells.Item(24, 10).Value = 0.0004555850957885808
$ws.Cells.Item(24, 13).Value = 1.018760333333333
$ws.Cells.Item(24, 14).Value = 3.056281
$ws.Cells.Item(24, 15).Value = 0.07126880433834142
$ws.Cells.Item(24, 16).Value = 0.07126880433834142
$ws.Cells.Item(24, 17).Value = 1.441543834146
$ws.Cells.Item(24, 18).Value = 12.973894507314
$ws.Cells.Item(24, 19).Value = 0.00003246900505122089
$ws.Cells.Item(24, 20).Value = 0.00003246900505122089

# Row 25
$ws.Cells.Item(25, 7).Value = 1.414998
$ws.Cells.Item(25, 8).Value = 4.244994
$ws.Cells.Item(25, 9).Value = 0.0004555850957885808
$ws.Cells.Item(25, 10).Value = 0.0004555850957885808
$ws.Cells.Item(25, 13).Value = 9.893154666666666
$ws.Cells.Item(25, 14).Value = 29.679464
$ws.Cells.Item(25, 15).Value = 0.6920894749804902
$ws.Cells.Item(25, 16).Value = 0.6920894749804903
$ws.Cells.Item(25, 17).Value = 13.998794067024
$ws.Cells.Item(25, 18).Value = 125.989146603216
$ws.Cells.Item(25, 19).Value = 0.0003153056497532552
$ws.Cells.Item(25, 20).Value = 0.0003153056497532553

# Row 26
$ws.Cells.Item(26, 7).Value = 1.414998
$ws.Cells.Item(26, 8).Value = 4.244994
$ws.Cells.Item(26, 9).Value = 0.0004555850957885808
$ws.Cells.Item(26, 10).Value = 0.0004555850957885808
$ws.Cells.Item(26, 13).Value = 0.1629343333333333
$ws.Cells.Item(26, 14).Value = 0.488803
$ws.Cells.Item(26, 15).Value = 0.01139829922935564
$ws.Cells.Item(26, 16).Value = 0.01139829922935564
$ws.Cells.Item(26, 17).Value = 0.2394788192984444
$ws.Cells.Item(26, 18).Value = 2.155309373686
$ws.Cells.Item(26, 19).Value = 0.000005393966391641539
$ws.Cells.Item(26, 20).Value = 0.000005393966391641539
